$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 6000
$ws.Range("N17").Value = -6336
$ws.Range("M17").ClearContents()
# Row 111
$ws.Range("H111").Value = 1131.5
$ws.Range("I111").Value = 1157.8667
$ws.Range("K111").Value = 3473.6001
$ws.Range("M111").Value = -406.6001000000001
# Row 113
$ws.Range("H113").Value = 1566.6666
$ws.Range("I113").Value = 1566.6666
$ws.Range("K113").Value = 1566.6666
$ws.Range("M113").Value = 1687.3334
# Row 135
$ws.Range("H135").Value = 739.8333
$ws.Range("I135").Value = 221.2
$ws.Range("J135").Value = 3333
$ws.Range("K135").Value = 1990.8
$ws.Range("L135").Value = 29997
$ws.Range("M135").Value = 544.2
$ws.Range("N135").Value = -35067
# Row 137
$ws.Range("H137").Value = 4683.1904
$ws.Range("I137").Value = 2156.5
$ws.Range("J137").Value = 6238.077
$ws.Range("K137").Value = 6469.5
$ws.Range("L137").Value = 18714.231
$ws.Range("M137").Value = -3919.5
$ws.Range("N137").Value = -23814.231
# Row 138
$ws.Range("H138").Value = 2393.7778
$ws.Range("I138").Value = 1270.25
$ws.Range("J138").Value = 3292.6
$ws.Range("K138").Value = 3810.75
$ws.Range("L138").Value = 9877.799999999999
$ws.Range("M138").Value = 1329.25
$ws.Range("N138").Value = -20157.8

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3133.5
$ws.Range("I2").Value = 3133.5
$ws.Range("K2").Value = 3133.5
$ws.Range("M2").Value = -3020.5
# Row 32
$ws.Range("H32").Value = 12106.387
$ws.Range("I32").Value = 9270.166999999999
$ws.Range("K32").Value = 9270.166999999999
$ws.Range("M32").Value = -8983.166999999999
# Row 74
$ws.Range("H74").Value = 2280.7727
$ws.Range("I74").Value = 1805.2354
$ws.Range("J74").Value = 3897.6
$ws.Range("K74").Value = 1805.2354
$ws.Range("L74").Value = 3897.6
$ws.Range("M74").Value = -931.2354
$ws.Range("N74").Value = -5645.6
# Row 77
$ws.Range("H77").Value = 2280.7727
$ws.Range("I77").Value = 1805.2354
$ws.Range("J77").Value = 3897.6
$ws.Range("K77").Value = 9026.177
$ws.Range("L77").Value = 19488
$ws.Range("M77").Value = -4658.177
$ws.Range("N77").Value = -28224
# Row 116
$ws.Range("H116").Value = 3133.5
$ws.Range("I116").Value = 3133.5
$ws.Range("K116").Value = 3133.5
$ws.Range("M116").Value = -839.5
# Row 122
$ws.Range("H122").Value = 974.5
$ws.Range("I122").Value = 974.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2923.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -473.5
$ws.Range("N122").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3133.5
$ws.Range("I3").Value = 3133.5
$ws.Range("K3").Value = 3133.5
$ws.Range("M3").Value = -3019.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 140.42857
$ws.Range("I7").Value = 113.833336
$ws.Range("K7").Value = 113.833336
$ws.Range("M7").Value = -0.8333360000000027
# Row 16
$ws.Range("H16").Value = 3429.8
$ws.Range("I16").Value = 2716.3333
$ws.Range("K16").Value = 2716.3333
$ws.Range("M16").Value = -2429.3333
# Row 31
$ws.Range("H31").Value = 5538.75
$ws.Range("I31").Value = 4582.773
$ws.Range("K31").Value = 4582.773
$ws.Range("M31").Value = -4287.773
# Row 34
$ws.Range("H34").Value = 5538.75
$ws.Range("I34").Value = 4582.773
$ws.Range("K34").Value = 4582.773
$ws.Range("M34").Value = -4380.773
# Row 58
$ws.Range("H58").Value = 2669.1428
$ws.Range("I58").Value = 1919.6666
$ws.Range("J58").Value = 7166
$ws.Range("K58").Value = 1919.6666
$ws.Range("L58").Value = 7166
$ws.Range("M58").Value = -1716.6666
$ws.Range("N58").Value = -7572
# Row 86
$ws.Range("H86").Value = 2787.5
$ws.Range("I86").Value = 2757.1428
$ws.Range("K86").Value = 2757.1428
$ws.Range("M86").Value = -1634.1428
# Row 89
$ws.Range("H89").Value = 2787.5
$ws.Range("I89").Value = 2757.1428
$ws.Range("K89").Value = 13785.714
$ws.Range("M89").Value = -8169.714
# Row 99
$ws.Range("H99").Value = 5026.75
$ws.Range("I99").Value = 4695.2666
$ws.Range("K99").Value = 4695.2666
$ws.Range("M99").Value = -3197.2666
# Row 107
$ws.Range("H107").Value = 386.70834
$ws.Range("I107").Value = 319.55
$ws.Range("J107").Value = 722.5
$ws.Range("K107").Value = 319.55
$ws.Range("L107").Value = 722.5
$ws.Range("M107").Value = 1600.45
$ws.Range("N107").Value = -4562.5
# Row 113
$ws.Range("H113").Value = 3429.8
$ws.Range("I113").Value = 2716.3333
$ws.Range("K113").Value = 2716.3333
$ws.Range("M113").Value = -546.3332999999998
# Row 126
$ws.Range("H126").Value = 5026.75
$ws.Range("I126").Value = 4695.2666
$ws.Range("K126").Value = 14085.7998
$ws.Range("M126").Value = -11615.7998
# Row 132
$ws.Range("H132").Value = 2664.2942
$ws.Range("I132").Value = 2206.125
$ws.Range("K132").Value = 6618.375
$ws.Range("M132").Value = -4088.375
# Row 136
$ws.Range("H136").Value = 2669.1428
$ws.Range("I136").Value = 1919.6666
$ws.Range("J136").Value = 7166
$ws.Range("K136").Value = 5758.9998
$ws.Range("L136").Value = 21498
$ws.Range("M136").Value = -3208.9998
$ws.Range("N136").Value = -26598

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 323.23077
$ws.Range("I11").Value = 600.5
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 1801.5
$ws.Range("L11").Value = 600
$ws.Range("M11").Value = -1661.5
$ws.Range("N11").Value = -880
# Row 23
$ws.Range("H23").Value = 122.38461
$ws.Range("J23").Value = 125.25
$ws.Range("L23").Value = 375.75
$ws.Range("N23").Value = -845.75
# Row 131
$ws.Range("H131").Value = 998.3333
$ws.Range("I131").Value = 998.3333
$ws.Range("J131").Value = 0
$ws.Range("K131").Value = 2994.9999
$ws.Range("L131").Value = 0
$ws.Range("M131").Value = 2045.0001
$ws.Range("N131").ClearContents()
# Row 132
$ws.Range("H132").Value = 3504.5908
$ws.Range("I132").Value = 4200.143
$ws.Range("J132").Value = 3180
$ws.Range("K132").Value = 37801.287
$ws.Range("L132").Value = 28620
$ws.Range("M132").Value = -35271.287
$ws.Range("N132").Value = -33680

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 93
$ws.Range("H93").Value = 47666.332
$ws.Range("J93").Value = 47666.332
$ws.Range("L93").Value = 47666.332
$ws.Range("N93").Value = -51410.332
# Row 117
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
# Row 122
$ws.Range("H122").Value = 5310.5
$ws.Range("I122").Value = 2844
$ws.Range("J122").Value = 7777
$ws.Range("K122").Value = 8532
$ws.Range("L122").Value = 23331
$ws.Range("M122").Value = -6082
$ws.Range("N122").Value = -28231

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 6659.2
$ws.Range("I136").Value = 5824
$ws.Range("K136").Value = 17472
$ws.Range("M136").Value = -14922

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3403.5173
$ws.Range("I122").Value = 2270.1875
$ws.Range("J122").Value = 4798.385
$ws.Range("K122").Value = 6810.5625
$ws.Range("L122").Value = 14395.155
$ws.Range("M122").Value = -4360.5625
$ws.Range("N122").Value = -19295.155
# Row 136
$ws.Range("H136").Value = 5386.2666
$ws.Range("I136").Value = 4028.3635
$ws.Range("J136").Value = 9120.5
$ws.Range("K136").Value = 12085.0905
$ws.Range("L136").Value = 27361.5
$ws.Range("M136").Value = -9535.0905
$ws.Range("N136").Value = -32461.5
